$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2949.6758
$ws.Range("J17").Value = 2949.6758
$ws.Range("L17").Value = 8849.027399999999
$ws.Range("N17").Value = -9185.027399999999
$ws.Range("H55").Value = 687.3333
$ws.Range("I55").Value = 174.16667
$ws.Range("J55").Value = 1713.6666
$ws.Range("K55").Value = 174.16667
$ws.Range("L55").Value = 1713.6666
$ws.Range("M55").Value = 39.83332999999999
$ws.Range("N55").Value = -2141.6666
$ws.Range("H80").Value = 465.9375
$ws.Range("I80").Value = 313.83334
$ws.Range("J80").Value = 557.2
$ws.Range("K80").Value = 941.5000200000001
$ws.Range("L80").Value = 1671.6
$ws.Range("M80").Value = 56.49997999999994
$ws.Range("N80").Value = -3667.6
$ws.Range("H83").Value = 465.9375
$ws.Range("I83").Value = 313.83334
$ws.Range("J83").Value = 557.2
$ws.Range("K83").Value = 2824.50006
$ws.Range("L83").Value = 5014.8
$ws.Range("M83").Value = 2167.49994
$ws.Range("N83").Value = -14998.8
$ws.Range("H112").Value = 1808.303
$ws.Range("I112").Value = 847.3333
$ws.Range("J112").Value = 1904.4
$ws.Range("K112").Value = 2541.9999
$ws.Range("L112").Value = 5713.200000000001
$ws.Range("M112").Value = -1433.9999
$ws.Range("N112").Value = -7929.200000000001
$ws.Range("H132").Value = 20835930
$ws.Range("I132").Value = 23812058
$ws.Range("J132").Value = 3046
$ws.Range("K132").Value = 71436174
$ws.Range("L132").Value = 9138
$ws.Range("M132").Value = -71433644
$ws.Range("N132").Value = -14198
$ws.Range("H135").Value = 2001.1555
$ws.Range("I135").Value = 1722.3489
$ws.Range("J135").Value = 7995.5
$ws.Range("K135").Value = 15501.1401
$ws.Range("L135").Value = 71959.5
$ws.Range("M135").Value = -12966.1401
$ws.Range("N135").Value = -77029.5
$ws.Range("H137").Value = 5238.3887
$ws.Range("I137").Value = 5286.1333
$ws.Range("J137").Value = 4999.6665
$ws.Range("K137").Value = 15858.3999
$ws.Range("L137").Value = 14998.9995
$ws.Range("M137").Value = -13308.3999
$ws.Range("N137").Value = -20098.9995
$ws.Range("H138").Value = 2502.3372
$ws.Range("I138").Value = 1214.3959
$ws.Range("J138").Value = 4129.2104
$ws.Range("K138").Value = 3643.1877
$ws.Range("L138").Value = 12387.6312
$ws.Range("M138").Value = 1496.8123
$ws.Range("N138").Value = -22667.6312

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5220.6313
$ws.Range("I2").Value = 384.66666
$ws.Range("K2").Value = 384.66666
$ws.Range("M2").Value = -271.66666
$ws.Range("H32").Value = 2683.7966
$ws.Range("I32").Value = 2437
$ws.Range("J32").Value = 16998
$ws.Range("K32").Value = 2437
$ws.Range("L32").Value = 16998
$ws.Range("M32").Value = -2150
$ws.Range("N32").Value = -17572
$ws.Range("H45").Value = 1605.5385
$ws.Range("I45").Value = 1696.1111
$ws.Range("J45").Value = 1401.75
$ws.Range("K45").Value = 1696.1111
$ws.Range("L45").Value = 1401.75
$ws.Range("M45").Value = -1319.1111
$ws.Range("N45").Value = -2155.75
$ws.Range("H61").Value = 2429.08
$ws.Range("I61").Value = 1687.6364
$ws.Range("J61").Value = 7866.3335
$ws.Range("K61").Value = 1687.6364
$ws.Range("L61").Value = 7866.3335
$ws.Range("M61").Value = -1475.6364
$ws.Range("N61").Value = -8290.333500000001
$ws.Range("H74").Value = 2723.75
$ws.Range("I74").Value = 1661.75
$ws.Range("J74").Value = 4847.75
$ws.Range("K74").Value = 1661.75
$ws.Range("L74").Value = 4847.75
$ws.Range("M74").Value = -787.75
$ws.Range("N74").Value = -6595.75
$ws.Range("H77").Value = 2723.75
$ws.Range("I77").Value = 1661.75
$ws.Range("J77").Value = 4847.75
$ws.Range("K77").Value = 8308.75
$ws.Range("L77").Value = 24238.75
$ws.Range("M77").Value = -3940.75
$ws.Range("N77").Value = -32974.75
$ws.Range("H102").Value = 31317992
$ws.Range("I102").Value = 38467530
$ws.Range("J102").Value = 336666.34
$ws.Range("K102").Value = 38467530
$ws.Range("L102").Value = 336666.34
$ws.Range("M102").Value = -38465908
$ws.Range("N102").Value = -339910.34
$ws.Range("H116").Value = 5220.6313
$ws.Range("I116").Value = 384.66666
$ws.Range("K116").Value = 384.66666
$ws.Range("M116").Value = 1909.33334
$ws.Range("H122").Value = 4561.1787
$ws.Range("I122").Value = 4774.4585
$ws.Range("K122").Value = 14323.3755
$ws.Range("M122").Value = -11873.3755
$ws.Range("H132").Value = 7210.9033
$ws.Range("I132").Value = 4518.3623
$ws.Range("K132").Value = 13555.0869
$ws.Range("M132").Value = -11025.0869
$ws.Range("H136").Value = 2429.08
$ws.Range("I136").Value = 1687.6364
$ws.Range("J136").Value = 7866.3335
$ws.Range("K136").Value = 5062.9092
$ws.Range("L136").Value = 23599.0005
$ws.Range("M136").Value = -2512.9092
$ws.Range("N136").Value = -28699.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5220.6313
$ws.Range("I3").Value = 384.66666
$ws.Range("K3").Value = 384.66666
$ws.Range("M3").Value = -270.66666
$ws.Range("H99").Value = 1771.4286
$ws.Range("I99").Value = 1636
$ws.Range("J99").Value = 2110
$ws.Range("K99").Value = 1636
$ws.Range("L99").Value = 2110
$ws.Range("M99").Value = -138
$ws.Range("N99").Value = -5106
$ws.Range("H134").Value = 4052.1555
$ws.Range("I134").Value = 1962.9688
$ws.Range("K134").Value = 5888.9064
$ws.Range("M134").Value = -3353.9064

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2411.8108
$ws.Range("I58").Value = 2433.6365
$ws.Range("J58").Value = 2231.75
$ws.Range("K58").Value = 2433.6365
$ws.Range("L58").Value = 2231.75
$ws.Range("M58").Value = -2230.6365
$ws.Range("N58").Value = -2637.75
$ws.Range("H105").Value = 9072.091
$ws.Range("I105").Value = 10732.723
$ws.Range("K105").Value = 10732.723
$ws.Range("M105").Value = -8985.723
$ws.Range("H132").Value = 1005.43335
$ws.Range("I132").Value = 869.96295
$ws.Range("J132").Value = 2224.6667
$ws.Range("K132").Value = 2609.88885
$ws.Range("L132").Value = 6674.000100000001
$ws.Range("M132").Value = -79.88884999999982
$ws.Range("N132").Value = -11734.0001
$ws.Range("H134").Value = 1164.8889
$ws.Range("I134").Value = 1150.5245
$ws.Range("J134").Value = 1244.5454
$ws.Range("K134").Value = 3451.5735
$ws.Range("L134").Value = 3733.6362
$ws.Range("M134").Value = -916.5735
$ws.Range("N134").Value = -8803.636200000001
$ws.Range("H136").Value = 2411.8108
$ws.Range("I136").Value = 2433.6365
$ws.Range("J136").Value = 2231.75
$ws.Range("K136").Value = 7300.9095
$ws.Range("L136").Value = 6695.25
$ws.Range("M136").Value = -4750.9095
$ws.Range("N136").Value = -11795.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15151941
$ws.Range("I2").Value = 185.75
$ws.Range("J2").Value = 33334048
$ws.Range("K2").Value = 1114.5
$ws.Range("L2").Value = 200004288
$ws.Range("M2").Value = -1001.5
$ws.Range("N2").Value = -200004514
$ws.Range("H38").Value = 358.0909
$ws.Range("I38").Value = 170.5
$ws.Range("J38").Value = 465.2857
$ws.Range("K38").Value = 511.5
$ws.Range("L38").Value = 1395.8571
$ws.Range("M38").Value = -164.5
$ws.Range("N38").Value = -2089.8571
$ws.Range("H40").Value = 286.77777
$ws.Range("I40").Value = 287.13333
$ws.Range("J40").Value = 285
$ws.Range("K40").Value = 1148.53332
$ws.Range("L40").Value = 1140
$ws.Range("M40").Value = -1079.53332
$ws.Range("N40").Value = -1278
$ws.Range("H56").Value = 18936.25
$ws.Range("I56").Value = 18936.25
$ws.Range("K56").Value = 18936.25
$ws.Range("M56").Value = -18406.25
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H117").Value = 80860.46000000001
$ws.Range("I117").Value = 742
$ws.Range("J117").Value = 104896
$ws.Range("K117").Value = 2226
$ws.Range("L117").Value = 314688
$ws.Range("M117").Value = 1216
$ws.Range("N117").Value = -321572
$ws.Range("H121").Value = 2982.5625
$ws.Range("I121").Value = 275
$ws.Range("J121").Value = 3369.3572
$ws.Range("K121").Value = 825
$ws.Range("L121").Value = 10108.0716
$ws.Range("M121").Value = 485
$ws.Range("N121").Value = -12728.0716
$ws.Range("H137").Value = 3222.625
$ws.Range("J137").Value = 3222.625
$ws.Range("L137").Value = 9667.875
$ws.Range("N137").Value = -19867.875
$ws.Range("H139").Value = 3820.9092
$ws.Range("I139").Value = 3203
$ws.Range("J139").Value = 10000
$ws.Range("K139").Value = 9609
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -4469
$ws.Range("N139").Value = -40280
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 537.6
$ws.Range("I97").Value = 517.5833
$ws.Range("J97").Value = 567.625
$ws.Range("K97").Value = 517.5833
$ws.Range("L97").Value = 567.625
$ws.Range("M97").Value = -21.58330000000001
$ws.Range("N97").Value = -1559.625
$ws.Range("H126").Value = 5064.3076
$ws.Range("I126").Value = 4433.6
$ws.Range("J126").Value = 7166.6665
$ws.Range("K126").Value = 13300.8
$ws.Range("L126").Value = 21499.9995
$ws.Range("M126").Value = -10830.8
$ws.Range("N126").Value = -26439.9995
$ws.Range("H132").Value = 893.02325
$ws.Range("I132").Value = 987.3143
$ws.Range("J132").Value = 480.5
$ws.Range("K132").Value = 2961.9429
$ws.Range("L132").Value = 1441.5
$ws.Range("M132").Value = -431.9429
$ws.Range("N132").Value = -6501.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3327.62
$ws.Range("I132").Value = 1975.2703
$ws.Range("K132").Value = 5925.810899999999
$ws.Range("M132").Value = -3395.810899999999
$ws.Range("H136").Value = 3928.4736
$ws.Range("I136").Value = 3046.4375
$ws.Range("K136").Value = 9139.3125
$ws.Range("M136").Value = -6589.3125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1809.8125
$ws.Range("I126").Value = 1803.6774
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5411.0322
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2941.0322
$ws.Range("N126").Value = -10940
$ws.Range("H136").Value = 4099.231
$ws.Range("I136").Value = 4360.087
$ws.Range("J136").Value = 2099.3333
$ws.Range("K136").Value = 13080.261
$ws.Range("L136").Value = 6297.999899999999
$ws.Range("M136").Value = -10530.261
$ws.Range("N136").Value = -11397.9999

